$d = $word.ActiveDocument

# --- 1. Drop the "_GoBack" bookmark from its old spot -----------------------
# In the original document, "_GoBack" sits right after the repository URL,
# at the very end of the "Repositorio: ..." paragraph. It needs to move to
# the end of the "Alfonso Sanchez ..." paragraph instead, so delete it here
# and re-create it later in the right place.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# --- 2. Add the "Blog:" label and URL to Alfonso's line ----------------------
# These become two separate runs, appended right after "U900428".
$rng = $d.Content
$rng.Find.Execute("U900428")
$rng.Collapse(0)
$rng.InsertAfter("  Blog: ")
$rng.Collapse(0)
$rng.InsertAfter("http://dsd2013.blogspot.com/")

# --- 3. Re-create "_GoBack" at the new end-of-paragraph position -------------
# A collapsed range placed exactly at "end of paragraph text" (i.e. right
# before the paragraph mark) trips this host's Bookmarks.Add into anchoring
# the bookmark somewhere else entirely. Work around it by first inserting a
# throwaway marker *after* the bookmark's intended position, anchoring the
# bookmark to the (now non-edge) spot right before the marker, and only then
# deleting the marker - normal text deletion keeps the bookmark correctly
# collapsed where it belongs.
$markRng = $d.Content
$markRng.Find.Execute("http://dsd2013.blogspot.com/")
$markRng.Collapse(0)
$markRng.InsertAfter("@@TMPMARK@@")

$anchorRng = $d.Content
$anchorRng.Find.Execute("http://dsd2013.blogspot.com/")
$anchorRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $anchorRng)

$cleanupRng = $d.Content
$cleanupRng.Find.Execute("@@TMPMARK@@")
$cleanupRng.Text = ""
